$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen the table with two new columns (P, Q) across every data row (1-115),
# mirroring Excel's behaviour of extending the used range when a blank cell
# is touched/formatted, without altering any existing styling.
$ws.Range("P1:Q115").Style = "Normal"

# Record the new source statement (S51 "stated in" -> Q-item) on every row
# that already carries a sourced claim, using the first empty column pair.
$ws.Range("D1").Value = "S51"
$ws.Range("E1").Value = "Q476"
$ws.Range("D2").Value = "S51"
$ws.Range("E2").Value = "Q476"
$ws.Range("D3").Value = "S51"
$ws.Range("E3").Value = "Q476"
$ws.Range("J4").Value = "S51"
$ws.Range("K4").Value = "Q476"
$ws.Range("D5").Value = "S51"
$ws.Range("E5").Value = "Q476"
$ws.Range("D6").Value = "S51"
$ws.Range("E6").Value = "Q476"
$ws.Range("F7").Value = "S51"
$ws.Range("G7").Value = "Q476"
$ws.Range("D15").Value = "S51"
$ws.Range("E15").Value = "Q475"
$ws.Range("D16").Value = "S51"
$ws.Range("E16").Value = "Q475"
$ws.Range("D17").Value = "S51"
$ws.Range("E17").Value = "Q475"
$ws.Range("D18").Value = "S51"
$ws.Range("E18").Value = "Q475"
$ws.Range("J19").Value = "S51"
$ws.Range("K19").Value = "Q475"
$ws.Range("D20").Value = "S51"
$ws.Range("E20").Value = "Q475"
$ws.Range("D21").Value = "S51"
$ws.Range("E21").Value = "Q475"
$ws.Range("F22").Value = "S51"
$ws.Range("G22").Value = "Q475"
$ws.Range("D30").Value = "S51"
$ws.Range("E30").Value = "Q474"
$ws.Range("D31").Value = "S51"
$ws.Range("E31").Value = "Q474"
$ws.Range("D32").Value = "S51"
$ws.Range("E32").Value = "Q474"
$ws.Range("D33").Value = "S51"
$ws.Range("E33").Value = "Q474"
$ws.Range("J34").Value = "S51"
$ws.Range("K34").Value = "Q474"
$ws.Range("D35").Value = "S51"
$ws.Range("E35").Value = "Q474"
$ws.Range("D36").Value = "S51"
$ws.Range("E36").Value = "Q474"
$ws.Range("F37").Value = "S51"
$ws.Range("G37").Value = "Q474"
$ws.Range("D48").Value = "S51"
$ws.Range("E48").Value = "Q400"
$ws.Range("D49").Value = "S51"
$ws.Range("E49").Value = "Q400"
$ws.Range("D50").Value = "S51"
$ws.Range("E50").Value = "Q400"
$ws.Range("D51").Value = "S51"
$ws.Range("E51").Value = "Q400"
$ws.Range("J52").Value = "S51"
$ws.Range("K52").Value = "Q400"
$ws.Range("D53").Value = "S51"
$ws.Range("E53").Value = "Q400"
$ws.Range("D54").Value = "S51"
$ws.Range("E54").Value = "Q400"
$ws.Range("F55").Value = "S51"
$ws.Range("G55").Value = "Q400"
$ws.Range("D64").Value = "S51"
$ws.Range("E64").Value = "Q381"
$ws.Range("D65").Value = "S51"
$ws.Range("E65").Value = "Q381"
$ws.Range("D66").Value = "S51"
$ws.Range("E66").Value = "Q381"
$ws.Range("D67").Value = "S51"
$ws.Range("E67").Value = "Q381"
$ws.Range("L68").Value = "S51"
$ws.Range("M68").Value = "Q381"
$ws.Range("D69").Value = "S51"
$ws.Range("E69").Value = "Q381"
$ws.Range("D70").Value = "S51"
$ws.Range("E70").Value = "Q381"
$ws.Range("F71").Value = "S51"
$ws.Range("G71").Value = "Q381"
$ws.Range("F72").Value = "S51"
$ws.Range("G72").Value = "Q381"
$ws.Range("F73").Value = "S51"
$ws.Range("G73").Value = "Q381"
$ws.Range("D91").Value = "S51"
$ws.Range("E91").Value = "Q470"
$ws.Range("D92").Value = "S51"
$ws.Range("E92").Value = "Q470"
$ws.Range("D93").Value = "S51"
$ws.Range("E93").Value = "Q470"
$ws.Range("D94").Value = "S51"
$ws.Range("E94").Value = "Q470"
$ws.Range("L95").Value = "S51"
$ws.Range("M95").Value = "Q470"
$ws.Range("D96").Value = "S51"
$ws.Range("E96").Value = "Q470"
$ws.Range("D97").Value = "S51"
$ws.Range("E97").Value = "Q470"
$ws.Range("F98").Value = "S51"
$ws.Range("G98").Value = "Q470"
$ws.Range("F99").Value = "S51"
$ws.Range("G99").Value = "Q470"
